# Adding Demographics Breakdown Chart
#
# Appends new "POPINC" (Population Increase) indicator rows to the bottom
# of the indicators table on Sheet1: one national (AUS) row followed by
# one row per state/territory, each tagged with Division "DEMOG" and
# Aggregation_Rule "Sum", with column D rebuilding the Mnemonic_Division
# key via the same "B&"_"&C" formula used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$mnemonics = @(
    "AUSPOPINC",
    "NSWPOPINC",
    "VICPOPINC",
    "QLDPOPINC",
    "ACTPOPINC",
    "TASPOPINC",
    "SAPOPINC",
    "NTPOPINC",
    "WAPOPINC"
)

$firstRow = 286
$lastRow = $firstRow + $mnemonics.Length - 1

# First row: set the Series_ID (column E) before the Mnemonic (column B)
# so "POPINC" lands in the shared-strings table right before "AUSPOPINC",
# matching how the indicator was originally keyed in off the series name.
$ws.Cells.Item($firstRow, 1).Value = "AID"
$ws.Cells.Item($firstRow, 5).Value = "POPINC"
$ws.Cells.Item($firstRow, 2).Value = $mnemonics[0]
$ws.Cells.Item($firstRow, 3).Value = "DEMOG"
$ws.Cells.Item($firstRow, 6).Value = "Sum"

# Remaining rows (the state/territory breakdown), filled in afterwards.
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "AID"
    $ws.Cells.Item($r, 3).Value = "DEMOG"
    $ws.Cells.Item($r, 5).Value = "POPINC"
    $ws.Cells.Item($r, 6).Value = "Sum"
}
for ($i = 1; $i -lt $mnemonics.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $mnemonics[$i]
}

# Column D: "Mnemonic_Division" helper formula, filled down across the
# whole new block first, then re-filled from the second row onward so the
# state rows form their own shared-formula group (same pattern already
# used for every other block of rows in this sheet).
$ws.Range("D" + $firstRow + ":D" + $lastRow).Formula = "=B" + $firstRow + "&`"_`"&C" + $firstRow
$ws.Range("D" + ($firstRow + 1) + ":D" + $lastRow).Formula = "=B" + ($firstRow + 1) + "&`"_`"&C" + ($firstRow + 1)

# Leave the selection where the user was working when the sheet was saved.
$ws.Range("G" + $firstRow + ":G" + $lastRow).Select() | Out-Null
